$wb = $excel.ActiveWorkbook

# --- Conversion sheet: insert a new column before "Reference" (old R) ---
# and label it "Interest rate". Everything to the right (old R..Y) shifts
# right by one column (new S..Z) automatically.
$wsConversion = $wb.Worksheets.Item("Conversion")
$wsConversion.Columns("R:R").Insert()
$wsConversion.Range("R1").Value = "Interest rate"

# --- Distribution sheet: add a new "Pumps" category row ---
$wsDistribution = $wb.Worksheets.Item("Distribution")
$wsDistribution.Range("A4").Value = "Pumps"

# --- Terminology sheet: restructure the cost breakdown ---
$wsTerminology = $wb.Worksheets.Item("Terminology")

# "Investment cost" category becomes "CAPEX"
$wsTerminology.Range("A2").Value = "CAPEX"

# Insert a new row for the "Financial cost" heading, pushing the old
# "Operation & Maintanence cost" block (rows 4-6) down to rows 5-7.
$wsTerminology.Rows("4:4").Insert()
$wsTerminology.Range("B4").Value = "Financial cost"

# The old "Operation& Maintanence cost" category (now on row 5) becomes "OPEX"
$wsTerminology.Range("A5").Value = "OPEX"

# The old "Fuel cost " row (now on row 6) becomes "Utility cost "
$wsTerminology.Range("B6").Value = "Utility cost "

# --- Make Distribution the active sheet/selection, matching the saved view ---
$wsDistribution.Activate()
$wsDistribution.Range("I1").Select()
